$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing ecoregion-101 data block (rows 2-16) into a new
# block for ecoregion 102 (rows 17-31), inserting the rows so nothing below
# is overwritten. This preserves the original (pre-tweak) values, including
# the not-yet-updated ProbMortality column.
$ws.Rows("2:16").Copy()
$ws.Rows("17:31").Insert()

# The newly inserted rows describe ecoregion 102, not 101.
$ws.Range("B17:B31").Value = 102

# Final tweak: bump ProbMortality (column E) for ecoregion 101 up to 0.95.
$ws.Range("E2:E16").Value = 0.95

# Reflect the edited range in the sheet's active selection.
$null = $ws.Range("E2:E16").Select()
